# ABET - Taller Diseno de Software.docx
# Commit: "explicacion en archivo de Word"
#
# 1) Append the LSP explanation paragraph text after "LSP:".
# 2) Re-home the "_GoBack" bookmark from the old "Una solucion..." edit
#    point to the new edit point inside the LSP text (Bookmarks.Add with
#    an existing bookmark name moves it, matching Word's own behaviour of
#    tracking "last edit location").
# 3) Normalise the "Una solucion..." paragraph back into a single run
#    (the bookmark that used to split it into two runs is gone now).
# 4) Fix up the footer's cached PAGE field result (now page 2 of 2).

$d = $word.ActiveDocument

# --- 1) & 2): LSP paragraph -------------------------------------------------

$lspTail = " En el método escuchar de la clase Carro, existe una condición que verifica si la radio está encendida, pero no hay una procedimiento que siga el programa en caso de que no lo esté. Por ello la solución  sería de implementar un método que se encargue de encender la radio, así el método escuchar verifica sin ningún problema el hecho de si está o no encendida la radio"
$lspEnd = " sin ningún problema."

$r = $d.Content
$r.Find.Execute("LSP:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter($lspTail)

# Remember where the bookmark needs to start (right after "...encendida la
# radio", before the final " sin ningún problema.").
$bookmarkPos = $r.End

$r.Collapse(0)
$r.InsertAfter($lspEnd)

# Re-create "_GoBack" at the new location -- since a bookmark name must be
# unique, adding it again simply relocates it away from its old spot in the
# "Una solución..." paragraph.
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- 3): merge the "Una solución..." runs back into one --------------------

$solucionText = "Una solución es implementar interfaces que se encarguen de los métodos que  se relaciones con todas la clases y cambiar el código para que la dependencia sea desde afuera para desacoplar cualquier clase de cualquier implementación de otro clase diferente."

$r2 = $d.Content
$r2.Find.Execute($solucionText, $true, $false, $false, $false, $false, $true, 1, $false, $solucionText, 2)

# --- 4): footer page number 1 -> 2 ------------------------------------------

$footer = $d.Sections(1).Footers(1)
$footerChar = $footer.Range.Characters(1)
$footerChar.Text = "2"
